$wb = $excel.ActiveWorkbook

# Remember which sheet was active before we start, so we can restore it
# (this workbook opens with "About" as the active/selected tab).
$originalActive = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("GDPbES")

# --- Rename existing "hydrogen" row (row 24) to "hydrogen combustion turbine" ---
$ws.Range("A24").Value = "hydrogen combustion turbine"
$ws.Range("A24").Font.Color = 0
$ws.Range("A24").VerticalAlignment = -4108

# --- Add a new row (25) for "hydrogen combined cycle" ---
$ws.Range("A25").Value = "hydrogen combined cycle"

# Copy A24's now-finalized formatting onto A25 rather than re-deriving it
# property-by-property (which would leave an orphaned intermediate style).
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill B25:AJ25 with 0, matching the pattern of the rest of the table.
for ($c = 2; $c -le 36; $c++) {
    $ws.Cells.Item(25, $c).Value = 0
}

# Update the on-screen selection to mirror the authored view state.
$ws.Activate()
$ws.Range("B25:AJ25").Select()

# Restore the originally active sheet/tab.
$originalActive.Activate()
